$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four distinct Neo4j/Cypher query strings live in B2:C4 (column B holds
# the per-tab "query" -- CasesTab/SamplesTab/FilesTab detail queries -- and
# column C holds the shared "StatQuery" count query, repeated for every row).
# All of them filter on sf.grouped_recurrence_score; update that recurrence
# score filter from "51-100" to "0-5" in every one of them.
$queryRange = $ws.Range("B2:C4")

foreach ($cell in $queryRange.Cells) {
    $text = $cell.Text
    if ($text -ne $null -and $text.Contains('["51-100"]')) {
        $cell.Value = $text.Replace('["51-100"]', '["0-5"]')
    }
}

# Move the active selection / scroll position: previously the sheet was
# scrolled to row 5 with C14 selected; now it's scrolled back up with E3
# selected as the active cell.
$ws.Range("E3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
